# Apply odds updates to Sheet1 of the workbook (rows 4 and 5).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 changes
$ws.Range("Q4").Value = 2.5
$ws.Range("R4").Value = 1.5
$ws.Range("V4").Value = 1.73

# Row 5 changes
$ws.Range("G5").Value = 1.7
$ws.Range("I5").Value = 5.75
$ws.Range("J5").Value = 2.38
$ws.Range("L5").Value = 6
$ws.Range("Q5").Value = 2.25
$ws.Range("R5").Value = 1.62
$ws.Range("U5").Value = 2.2
$ws.Range("V5").Value = 1.62
$ws.Range("Z5").Value = 12
$ws.Range("AE5").Value = 21
$ws.Range("AK5").Value = 67
$ws.Range("AL5").Value = 51
$ws.Range("AO5").Value = 9
$ws.Range("AV5").Value = 81
$ws.Range("AW5").Value = 7
$ws.Range("AX5").Value = 34
